$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.151.08"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.37%  '

$ws.Range('D3').Value = "'1.878.49"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.76%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = "'243.61"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.51%  '

$ws.Range('D6').Value = "'0.9986"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.06%  '

$ws.Range('D7').Value = "'0.4914"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.16%  '

$ws.Range('D8').Value = "'0.2920"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.99%  '

$ws.Range('D9').Value = "'0.06611"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.41%  '

$ws.Range('D10').Value = "'1.880.27"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.67%  '

$ws.Range('D11').Value = "'16.55"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.54%  '

$ws.Range('D12').Value = "'0.07191"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.55%  '

$ws.Range('D13').Value = "'0.6681"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.50%  '

$ws.Range('D14').Value = "'86.36"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.35%  '

$ws.Range('D15').Value = "'4.928"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.62%  '

$ws.Range('D16').Value = "'30.106.57"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.30%  '

$ws.Range('D17').Value = "'0.000007828"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.21%  '

$ws.Range('E18').Value = '  -0.10%  '

$ws.Range('D19').Value = "'12.84"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.46%  '

$ws.Range('D20').Value = "'2.125.79"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.66%  '

$ws.Range('D21').Value = "'0.9982"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.01%  '

$ws.Range('D22').Value = "'4.794"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.63%  '

$ws.Range('D23').Value = "'5.855"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.20%  '

$ws.Range('D24').Value = "'9.170"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.47%  '

$ws.Range('D25').Value = "'153.27"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.77%  '

$ws.Range('D26').Value = "'145.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.08%  '

$ws.Range('D27').Value = "'17.00"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.51%  '

$ws.Range('D28').Value = "'1.904"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.88%  '

$ws.Range('D29').Value = "'1.387"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.94%  '

$ws.Range('D30').Value = "'4.218"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.07%  '

$ws.Range('D31').Value = "'0.08797"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.50%  '

$ws.Range('D32').Value = "'4.007"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.70%  '

$ws.Range('D33').Value = "'0.05083"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.10%  '

$ws.Range('D34').Value = "'0.7213"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.19%  '

$ws.Range('D35').Value = "'1.117"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.16%  '

$ws.Range('E36').Value = '  -0.42%  '

$ws.Range('D37').Value = "'0.01845"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +11.41%  '

$ws.Range('D38').Value = "'2.683"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.45%  '

$ws.Range('D39').Value = "'2.168"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.84%  '

$ws.Range('D40').Value = "'0.9317"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.75%  '

$ws.Range('D41').Value = "'5.791"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.72%  '

$ws.Range('D42').Value = "'0.4239"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.66%  '

$ws.Range('D43').Value = "'0.9982"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.03%  '

$ws.Range('D44').Value = "'103.41"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.04%  '

$ws.Range('D45').Value = "'7.409"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.94%  '

$ws.Range('E46').Value = '  +1.89%  '

$ws.Range('D47').Value = "'0.05692"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.70%  '

$ws.Range('D48').Value = "'32.91"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.35%  '

$ws.Range('D49').Value = "'8.330"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.37%  '

$ws.Range('D50').Value = "'0.3771"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.07%  '

$ws.Range('D51').Value = "'1.345"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.23%  '
